$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 11:52"

# Rewrite the data table (rows 4-69), now re-sorted descending by "Casos totales" (col B)
# after incorporating the day-over-day updated figures.
$ws.Cells.Item(4, 1).Value = "Madrid"
$ws.Cells.Item(4, 2).Value = 60487
$ws.Cells.Item(4, 3).Value = 33645
$ws.Cells.Item(4, 4).Value = 19158
$ws.Cells.Item(4, 5).Value = 7684
$ws.Cells.Item(5, 1).Value = "Cataluña"
$ws.Cells.Item(5, 2).Value = 44892
$ws.Cells.Item(5, 3).Value = 15477
$ws.Cells.Item(5, 4).Value = 25072
$ws.Cells.Item(5, 5).Value = 4343
$ws.Cells.Item(6, 1).Value = "Castilla-La Mancha"
$ws.Cells.Item(6, 2).Value = 17557
$ws.Cells.Item(6, 3).Value = 4577
$ws.Cells.Item(6, 4).Value = 10792
$ws.Cells.Item(6, 5).Value = 2188
$ws.Cells.Item(7, 1).Value = "Castilla y Leon"
$ws.Cells.Item(7, 2).Value = 17402
$ws.Cells.Item(7, 3).Value = 5777
$ws.Cells.Item(7, 4).Value = 10043
$ws.Cells.Item(7, 5).Value = 1582
$ws.Cells.Item(8, 1).Value = "Pais Vasco"
$ws.Cells.Item(8, 2).Value = 13436
$ws.Cells.Item(8, 3).Value = 8136
$ws.Cells.Item(8, 4).Value = 4133
$ws.Cells.Item(8, 5).Value = 1167
$ws.Cells.Item(9, 1).Value = "Andalucia"
$ws.Cells.Item(9, 2).Value = 11845
$ws.Cells.Item(9, 3).Value = 3847
$ws.Cells.Item(9, 4).Value = 6919
$ws.Cells.Item(9, 5).Value = 1079
$ws.Cells.Item(10, 1).Value = "Galicia"
$ws.Cells.Item(10, 2).Value = 8805
$ws.Cells.Item(10, 3).Value = 1676
$ws.Cells.Item(10, 4).Value = 6754
$ws.Cells.Item(10, 5).Value = 375
$ws.Cells.Item(11, 1).Value = "Bizkaia/Vizcaya"
$ws.Cells.Item(11, 2).Value = 7045
$ws.Cells.Item(11, 3).Value = 7124
$ws.Cells.Item(11, 4).Value = 4423
$ws.Cells.Item(11, 5).Value = 551
$ws.Cells.Item(12, 1).Value = "Ciudad Real"
$ws.Cells.Item(12, 2).Value = 6358
$ws.Cells.Item(12, 3).Value = 4178
$ws.Cells.Item(12, 4).Value = 10597
$ws.Cells.Item(12, 5).Value = 802
$ws.Cells.Item(13, 1).Value = "Aragon"
$ws.Cells.Item(13, 2).Value = 5237
$ws.Cells.Item(13, 3).Value = 1680
$ws.Cells.Item(13, 4).Value = 2876
$ws.Cells.Item(13, 5).Value = 681
$ws.Cells.Item(14, 1).Value = "Valencia/Valencia"
$ws.Cells.Item(14, 2).Value = 5131
$ws.Cells.Item(14, 3).Value = 2194
$ws.Cells.Item(14, 4).Value = 2583
$ws.Cells.Item(14, 5).Value = 515
$ws.Cells.Item(15, 1).Value = "Navarra"
$ws.Cells.Item(15, 2).Value = 4994
$ws.Cells.Item(15, 3).Value = 1445
$ws.Cells.Item(15, 4).Value = 3140
$ws.Cells.Item(15, 5).Value = 409
$ws.Cells.Item(16, 1).Value = "Toledo"
$ws.Cells.Item(16, 2).Value = 3938
$ws.Cells.Item(16, 3).Value = 4178
$ws.Cells.Item(16, 4).Value = 10597
$ws.Cells.Item(16, 5).Value = 504
$ws.Cells.Item(17, 1).Value = "La Rioja"
$ws.Cells.Item(17, 2).Value = 3839
$ws.Cells.Item(17, 3).Value = 1937
$ws.Cells.Item(17, 4).Value = 1599
$ws.Cells.Item(17, 5).Value = 303
$ws.Cells.Item(18, 1).Value = "Albacete"
$ws.Cells.Item(18, 2).Value = 3754
$ws.Cells.Item(18, 3).Value = 4178
$ws.Cells.Item(18, 4).Value = 10597
$ws.Cells.Item(18, 5).Value = 373
$ws.Cells.Item(19, 1).Value = "Zaragoza"
$ws.Cells.Item(19, 2).Value = 3678
$ws.Cells.Item(19, 3).Value = 913
$ws.Cells.Item(19, 4).Value = 2274
$ws.Cells.Item(19, 5).Value = 491
$ws.Cells.Item(20, 1).Value = "Alacant/Alicante"
$ws.Cells.Item(20, 2).Value = 3476
$ws.Cells.Item(20, 3).Value = 1677
$ws.Cells.Item(20, 4).Value = 1443
$ws.Cells.Item(20, 5).Value = 388
$ws.Cells.Item(21, 1).Value = "Extremadura"
$ws.Cells.Item(21, 2).Value = 3260
$ws.Cells.Item(21, 3).Value = 1237
$ws.Cells.Item(21, 4).Value = 1616
$ws.Cells.Item(21, 5).Value = 407
$ws.Cells.Item(22, 1).Value = "Araba/Alava"
$ws.Cells.Item(22, 2).Value = 3241
$ws.Cells.Item(22, 3).Value = 7124
$ws.Cells.Item(22, 4).Value = 4423
$ws.Cells.Item(22, 5).Value = 318
$ws.Cells.Item(23, 1).Value = "Valladolid"
$ws.Cells.Item(23, 2).Value = 3154
$ws.Cells.Item(23, 3).Value = 1070
$ws.Cells.Item(23, 4).Value = 1824
$ws.Cells.Item(23, 5).Value = 260
$ws.Cells.Item(24, 1).Value = "Salamanca"
$ws.Cells.Item(24, 2).Value = 2602
$ws.Cells.Item(24, 3).Value = 794
$ws.Cells.Item(24, 4).Value = 1521
$ws.Cells.Item(24, 5).Value = 287
$ws.Cells.Item(25, 1).Value = "Malaga"
$ws.Cells.Item(25, 2).Value = 2531
$ws.Cells.Item(25, 3).Value = 869
$ws.Cells.Item(25, 4).Value = 1439
$ws.Cells.Item(25, 5).Value = 223
$ws.Cells.Item(26, 1).Value = "Asturias"
$ws.Cells.Item(26, 2).Value = 2453
$ws.Cells.Item(26, 3).Value = 665
$ws.Cells.Item(26, 4).Value = 1565
$ws.Cells.Item(26, 5).Value = 223
$ws.Cells.Item(27, 1).Value = "Segovia"
$ws.Cells.Item(27, 2).Value = 2406
$ws.Cells.Item(27, 3).Value = 656
$ws.Cells.Item(27, 4).Value = 1578
$ws.Cells.Item(27, 5).Value = 172
$ws.Cells.Item(28, 1).Value = "Leon"
$ws.Cells.Item(28, 2).Value = 2403
$ws.Cells.Item(28, 3).Value = 1076
$ws.Cells.Item(28, 4).Value = 1024
$ws.Cells.Item(28, 5).Value = 303
$ws.Cells.Item(29, 1).Value = "Gipuzkoa/Guipuzcoa"
$ws.Cells.Item(29, 2).Value = 2342
$ws.Cells.Item(29, 3).Value = 7124
$ws.Cells.Item(29, 4).Value = 4423
$ws.Cells.Item(29, 5).Value = 212
$ws.Cells.Item(30, 1).Value = "Sevilla"
$ws.Cells.Item(30, 2).Value = 2329
$ws.Cells.Item(30, 3).Value = 459
$ws.Cells.Item(30, 4).Value = 1658
$ws.Cells.Item(30, 5).Value = 212
$ws.Cells.Item(31, 1).Value = "Caceres"
$ws.Cells.Item(31, 2).Value = 2220
$ws.Cells.Item(31, 3).Value = 422
$ws.Cells.Item(31, 4).Value = 1482
$ws.Cells.Item(31, 5).Value = 316
$ws.Cells.Item(32, 1).Value = "Cantabria"
$ws.Cells.Item(32, 2).Value = 2206
$ws.Cells.Item(32, 3).Value = 850
$ws.Cells.Item(32, 4).Value = 1183
$ws.Cells.Item(32, 5).Value = 173
$ws.Cells.Item(33, 1).Value = "Gran Canaria"
$ws.Cells.Item(33, 2).Value = 2113
$ws.Cells.Item(33, 3).Value = 969
$ws.Cells.Item(33, 4).Value = 1023
$ws.Cells.Item(33, 5).Value = 121
$ws.Cells.Item(34, 1).Value = "Granada"
$ws.Cells.Item(34, 2).Value = 2023
$ws.Cells.Item(34, 3).Value = 616
$ws.Cells.Item(34, 4).Value = 1211
$ws.Cells.Item(34, 5).Value = 196
$ws.Cells.Item(35, 1).Value = "A Coruña"
$ws.Cells.Item(35, 2).Value = 1969
$ws.Cells.Item(35, 3).Value = 333
$ws.Cells.Item(35, 4).Value = 1788
$ws.Cells.Item(35, 5).Value = 67
$ws.Cells.Item(36, 1).Value = "Murcia"
$ws.Cells.Item(36, 2).Value = 1721
$ws.Cells.Item(36, 3).Value = 766
$ws.Cells.Item(36, 4).Value = 831
$ws.Cells.Item(36, 5).Value = 124
$ws.Cells.Item(37, 1).Value = "Burgos"
$ws.Cells.Item(37, 2).Value = 1567
$ws.Cells.Item(37, 3).Value = 642
$ws.Cells.Item(37, 4).Value = 757
$ws.Cells.Item(37, 5).Value = 168
$ws.Cells.Item(38, 1).Value = "Pontevedra"
$ws.Cells.Item(38, 2).Value = 1536
$ws.Cells.Item(38, 3).Value = 333
$ws.Cells.Item(38, 4).Value = 1411
$ws.Cells.Item(38, 5).Value = 30
$ws.Cells.Item(39, 1).Value = "Guadalajara"
$ws.Cells.Item(39, 2).Value = 1431
$ws.Cells.Item(39, 3).Value = 4178
$ws.Cells.Item(39, 4).Value = 10597
$ws.Cells.Item(39, 5).Value = 186
$ws.Cells.Item(40, 1).Value = "Tenerife"
$ws.Cells.Item(40, 2).Value = 1391
$ws.Cells.Item(40, 3).Value = 813
$ws.Cells.Item(40, 4).Value = 798
$ws.Cells.Item(40, 5).Value = 119
$ws.Cells.Item(41, 1).Value = "Cuenca"
$ws.Cells.Item(41, 2).Value = 1315
$ws.Cells.Item(41, 3).Value = 4178
$ws.Cells.Item(41, 4).Value = 10597
$ws.Cells.Item(41, 5).Value = 156
$ws.Cells.Item(42, 1).Value = "Jaen"
$ws.Cells.Item(42, 2).Value = 1297
$ws.Cells.Item(42, 3).Value = 253
$ws.Cells.Item(42, 4).Value = 907
$ws.Cells.Item(42, 5).Value = 137
$ws.Cells.Item(43, 1).Value = "Cordoba"
$ws.Cells.Item(43, 2).Value = 1271
$ws.Cells.Item(43, 3).Value = 371
$ws.Cells.Item(43, 4).Value = 826
$ws.Cells.Item(43, 5).Value = 74
$ws.Cells.Item(44, 1).Value = "Castello/Castellon"
$ws.Cells.Item(44, 2).Value = 1257
$ws.Cells.Item(44, 3).Value = 435
$ws.Cells.Item(44, 4).Value = 739
$ws.Cells.Item(44, 5).Value = 139
$ws.Cells.Item(45, 1).Value = "Soria"
$ws.Cells.Item(45, 2).Value = 1243
$ws.Cells.Item(45, 3).Value = 299
$ws.Cells.Item(45, 4).Value = 848
$ws.Cells.Item(45, 5).Value = 96
$ws.Cells.Item(46, 1).Value = "Avila"
$ws.Cells.Item(46, 2).Value = 1155
$ws.Cells.Item(46, 3).Value = 450
$ws.Cells.Item(46, 4).Value = 596
$ws.Cells.Item(46, 5).Value = 109
$ws.Cells.Item(47, 1).Value = "Cadiz"
$ws.Cells.Item(47, 2).Value = 1139
$ws.Cells.Item(47, 3).Value = 283
$ws.Cells.Item(47, 4).Value = 781
$ws.Cells.Item(47, 5).Value = 75
$ws.Cells.Item(48, 1).Value = "Badajoz"
$ws.Cells.Item(48, 2).Value = 1023
$ws.Cells.Item(48, 3).Value = 436
$ws.Cells.Item(48, 4).Value = 514
$ws.Cells.Item(48, 5).Value = 73
$ws.Cells.Item(49, 1).Value = "Ourense"
$ws.Cells.Item(49, 2).Value = 751
$ws.Cells.Item(49, 3).Value = 333
$ws.Cells.Item(49, 4).Value = 660
$ws.Cells.Item(49, 5).Value = 22
$ws.Cells.Item(50, 1).Value = "Palencia"
$ws.Cells.Item(50, 2).Value = 716
$ws.Cells.Item(50, 3).Value = 223
$ws.Cells.Item(50, 4).Value = 432
$ws.Cells.Item(50, 5).Value = 61
$ws.Cells.Item(51, 1).Value = "Zamora"
$ws.Cells.Item(51, 2).Value = 611
$ws.Cells.Item(51, 3).Value = 219
$ws.Cells.Item(51, 4).Value = 327
$ws.Cells.Item(51, 5).Value = 65
$ws.Cells.Item(52, 1).Value = "Huesca"
$ws.Cells.Item(52, 2).Value = 601
$ws.Cells.Item(52, 3).Value = 137
$ws.Cells.Item(52, 4).Value = 384
$ws.Cells.Item(52, 5).Value = 80
$ws.Cells.Item(53, 1).Value = "Lugo"
$ws.Cells.Item(53, 2).Value = 586
$ws.Cells.Item(53, 3).Value = 333
$ws.Cells.Item(53, 4).Value = 520
$ws.Cells.Item(53, 5).Value = 11
$ws.Cells.Item(54, 1).Value = "Teruel"
$ws.Cells.Item(54, 2).Value = 541
$ws.Cells.Item(54, 3).Value = 117
$ws.Cells.Item(54, 4).Value = 359
$ws.Cells.Item(54, 5).Value = 65
$ws.Cells.Item(55, 1).Value = "Almeria"
$ws.Cells.Item(55, 2).Value = 458
$ws.Cells.Item(55, 3).Value = 147
$ws.Cells.Item(55, 4).Value = 269
$ws.Cells.Item(55, 5).Value = 42
$ws.Cells.Item(56, 1).Value = "Huelva"
$ws.Cells.Item(56, 2).Value = 377
$ws.Cells.Item(56, 3).Value = 118
$ws.Cells.Item(56, 4).Value = 225
$ws.Cells.Item(56, 5).Value = 34
$ws.Cells.Item(57, 1).Value = "Mallorca"
$ws.Cells.Item(57, 2).Value = 210
$ws.Cells.Item(57, 3).Value = 18
$ws.Cells.Item(57, 4).Value = 194
$ws.Cells.Item(57, 5).Value = 12
$ws.Cells.Item(58, 1).Value = "Ceuta"
$ws.Cells.Item(58, 2).Value = 118
$ws.Cells.Item(58, 3).Value = 76
$ws.Cells.Item(58, 4).Value = 38
$ws.Cells.Item(58, 5).Value = 4
$ws.Cells.Item(59, 1).Value = "Melilla"
$ws.Cells.Item(59, 2).Value = 106
$ws.Cells.Item(59, 3).Value = 55
$ws.Cells.Item(59, 4).Value = 49
$ws.Cells.Item(59, 5).Value = 2
$ws.Cells.Item(60, 1).Value = "La Palma"
$ws.Cells.Item(60, 2).Value = 74
$ws.Cells.Item(60, 3).Value = 25
$ws.Cells.Item(60, 4).Value = 46
$ws.Cells.Item(60, 5).Value = 3
$ws.Cells.Item(61, 1).Value = "Lanzarote"
$ws.Cells.Item(61, 2).Value = 68
$ws.Cells.Item(61, 3).Value = 17
$ws.Cells.Item(61, 4).Value = 49
$ws.Cells.Item(61, 5).Value = 2
$ws.Cells.Item(62, 1).Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Cells.Item(62, 2).Value = 58
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 58
$ws.Cells.Item(62, 5).Value = 3
$ws.Cells.Item(63, 1).Value = "Fuerteventura"
$ws.Cells.Item(63, 2).Value = 24
$ws.Cells.Item(63, 3).Value = 18
$ws.Cells.Item(63, 4).Value = 6
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 1).Value = "Ibiza"
$ws.Cells.Item(64, 2).Value = 21
$ws.Cells.Item(64, 3).Value = 18
$ws.Cells.Item(64, 4).Value = 20
$ws.Cells.Item(64, 5).Value = 1
$ws.Cells.Item(65, 1).Value = "Menorca"
$ws.Cells.Item(65, 2).Value = 15
$ws.Cells.Item(65, 3).Value = 18
$ws.Cells.Item(65, 4).Value = 13
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 1).Value = "Arroyo de la Luz"
$ws.Cells.Item(66, 2).Value = 7
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 4).Value = 7
$ws.Cells.Item(66, 5).Value = 0
$ws.Cells.Item(67, 1).Value = "La Gomera"
$ws.Cells.Item(67, 2).Value = 7
$ws.Cells.Item(67, 3).Value = 5
$ws.Cells.Item(67, 4).Value = 2
$ws.Cells.Item(67, 5).Value = 0
$ws.Cells.Item(68, 1).Value = "El Hierro"
$ws.Cells.Item(68, 2).Value = 1
$ws.Cells.Item(68, 3).Value = 1
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 1).Value = "Formentera"
$ws.Cells.Item(69, 2).Value = 0
$ws.Cells.Item(69, 3).Value = 10
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 8
